# Auto-generated PowerShell-style Excel COM-interop script
# Applies updated market-data values (H..N columns) to matching rows
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(138, 8).Value = 18194726  # H138: 18193342 -> 18194726
$ws.Cells.Item(138, 9).Value = 50004600  # I138: 50004000 -> 50004600
$ws.Cells.Item(138, 10).Value = 17655  # J138: 15824 -> 17655
$ws.Cells.Item(138, 11).Value = 150013800  # K138: 150012000 -> 150013800
$ws.Cells.Item(138, 12).Value = 52965  # L138: 47472 -> 52965
$ws.Cells.Item(138, 13).Value = -150008660  # M138: -150006860 -> -150008660
$ws.Cells.Item(138, 14).Value = -63245  # N138: -57752 -> -63245

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 33282.242  # H32: 34706.82 -> 33282.242
$ws.Cells.Item(32, 9).Value = 27374.521  # I32: 27333.217 -> 27374.521
$ws.Cells.Item(32, 10).Value = 67251.625  # J32: 119503.25 -> 67251.625
$ws.Cells.Item(32, 11).Value = 27374.521  # K32: 27333.217 -> 27374.521
$ws.Cells.Item(32, 12).Value = 67251.625  # L32: 119503.25 -> 67251.625
$ws.Cells.Item(32, 13).Value = -27087.521  # M32: -27046.217 -> -27087.521
$ws.Cells.Item(32, 14).Value = -67825.625  # N32: -120077.25 -> -67825.625
$ws.Cells.Item(138, 8).Value = 70000  # H138: 64221.25 -> 70000
$ws.Cells.Item(138, 10).Value = 70000  # J138: 64221.25 -> 70000
$ws.Cells.Item(138, 12).Value = 70000  # L138: 64221.25 -> 70000
$ws.Cells.Item(138, 14).Value = -80280  # N138: -74501.25 -> -80280

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(132, 8).Value = 72642.45  # H132: 70717.125 -> 72642.45
$ws.Cells.Item(132, 10).Value = 72642.45  # J132: 70717.125 -> 72642.45
$ws.Cells.Item(132, 12).Value = 72642.45  # L132: 70717.125 -> 72642.45
$ws.Cells.Item(132, 14).Value = -82762.45  # N132: -80837.125 -> -82762.45
$ws.Cells.Item(140, 8).Value = 102260  # H140: 93593.336 -> 102260
$ws.Cells.Item(140, 10).Value = 102260  # J140: 93593.336 -> 102260
$ws.Cells.Item(140, 12).Value = 102260  # L140: 93593.336 -> 102260
$ws.Cells.Item(140, 14).Value = -112620  # N140: -103953.336 -> -112620

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 5192.8  # H31: 5035.316 -> 5192.8
$ws.Cells.Item(31, 9).Value = 4753.5  # I31: 4370.5557 -> 4753.5
$ws.Cells.Item(31, 10).Value = 5632.1  # J31: 5633.6 -> 5632.1
$ws.Cells.Item(31, 11).Value = 4753.5  # K31: 4370.5557 -> 4753.5
$ws.Cells.Item(31, 12).Value = 5632.1  # L31: 5633.6 -> 5632.1
$ws.Cells.Item(31, 13).Value = -4458.5  # M31: -4075.5557 -> -4458.5
$ws.Cells.Item(31, 14).Value = -6222.1  # N31: -6223.6 -> -6222.1
$ws.Cells.Item(34, 8).Value = 5192.8  # H34: 5035.316 -> 5192.8
$ws.Cells.Item(34, 9).Value = 4753.5  # I34: 4370.5557 -> 4753.5
$ws.Cells.Item(34, 10).Value = 5632.1  # J34: 5633.6 -> 5632.1
$ws.Cells.Item(34, 11).Value = 4753.5  # K34: 4370.5557 -> 4753.5
$ws.Cells.Item(34, 12).Value = 5632.1  # L34: 5633.6 -> 5632.1
$ws.Cells.Item(34, 13).Value = -4551.5  # M34: -4168.5557 -> -4551.5
$ws.Cells.Item(34, 14).Value = -6036.1  # N34: -6037.6 -> -6036.1
$ws.Cells.Item(107, 8).Value = 570.25806  # H107: 611.85187 -> 570.25806
$ws.Cells.Item(107, 9).Value = 411.73685  # I107: 444.33334 -> 411.73685
$ws.Cells.Item(107, 11).Value = 411.73685  # K107: 444.33334 -> 411.73685
$ws.Cells.Item(107, 13).Value = 1508.26315  # M107: 1475.66666 -> 1508.26315
$ws.Cells.Item(133, 8).Value = 55326  # H133: 49326 -> 55326
$ws.Cells.Item(133, 10).Value = 55326  # J133: 49326 -> 55326
$ws.Cells.Item(133, 12).Value = 55326  # L133: 49326 -> 55326
$ws.Cells.Item(133, 14).Value = -60386  # N133: -54386 -> -60386
$ws.Cells.Item(135, 8).Value = 139784  # H135: 199866.67 -> 139784
$ws.Cells.Item(135, 10).Value = 139784  # J135: 199866.67 -> 139784
$ws.Cells.Item(135, 12).Value = 139784  # L135: 199866.67 -> 139784
$ws.Cells.Item(135, 14).Value = -149924  # N135: -210006.67 -> -149924

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(113, 8).Value = 14262.375  # H113: 14337.375 -> 14262.375
$ws.Cells.Item(113, 9).Value = 700  # I113: 900 -> 700
$ws.Cells.Item(113, 11).Value = 2100  # K113: 2700 -> 2100
$ws.Cells.Item(113, 13).Value = 70  # M113: -530 -> 70
$ws.Cells.Item(129, 8).Value = 2382781  # H129: 1853457.2 -> 2382781
$ws.Cells.Item(129, 9).Value = 750  # I129: 590 -> 750
$ws.Cells.Item(129, 10).Value = 2633521  # J129: 2175695 -> 2633521
$ws.Cells.Item(129, 11).Value = 2250  # K129: 1770 -> 2250
$ws.Cells.Item(129, 12).Value = 7900563  # L129: 6527085 -> 7900563
$ws.Cells.Item(129, 13).Value = 2750  # M129: 3230 -> 2750
$ws.Cells.Item(129, 14).Value = -7910563  # N129: -6537085 -> -7910563
$ws.Cells.Item(131, 8).Value = 21742846  # H131: 22225998 -> 21742846
$ws.Cells.Item(131, 10).Value = 22728866  # J131: 23257422 -> 22728866
$ws.Cells.Item(131, 12).Value = 68186598  # L131: 69772266 -> 68186598
$ws.Cells.Item(131, 14).Value = -68196678  # N131: -69782346 -> -68196678
$ws.Cells.Item(133, 8).Value = 2554.2144  # H133: 2802.1428 -> 2554.2144
$ws.Cells.Item(133, 9).Value = 1502.375  # I133: 1641.4286 -> 1502.375
$ws.Cells.Item(133, 10).Value = 3956.6667  # J133: 3962.8572 -> 3956.6667
$ws.Cells.Item(133, 11).Value = 4507.125  # K133: 4924.2858 -> 4507.125
$ws.Cells.Item(133, 12).Value = 11870.0001  # L133: 11888.5716 -> 11870.0001
$ws.Cells.Item(133, 13).Value = 552.875  # M133: 135.7142000000003 -> 552.875
$ws.Cells.Item(133, 14).Value = -21990.0001  # N133: -22008.5716 -> -21990.0001
$ws.Cells.Item(139, 8).Value = 1946.7241  # H139: 1994.8214 -> 1946.7241
$ws.Cells.Item(139, 9).Value = 1589.9524  # I139: 1639.45 -> 1589.9524
$ws.Cells.Item(139, 11).Value = 4769.857199999999  # K139: 4918.35 -> 4769.857199999999
$ws.Cells.Item(139, 13).Value = 370.1428000000005  # M139: 221.6499999999996 -> 370.1428000000005
$ws.Cells.Item(140, 8).Value = 6328.769  # H140: 3722.889 -> 6328.769
$ws.Cells.Item(140, 9).Value = 3602  # I140: 2491.4285 -> 3602
$ws.Cells.Item(140, 11).Value = 10806  # K140: 7474.2855 -> 10806
$ws.Cells.Item(140, 13).Value = -5626  # M140: -2294.2855 -> -5626

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(93, 8).Value = 0  # H93: 20250 -> 0
$ws.Cells.Item(93, 10).Value = 0  # J93: 20250 -> 0
$ws.Cells.Item(93, 12).Value = 0  # L93: 20250 -> 0
$ws.Cells.Item(93, 14).ClearContents()  # N93: -23994 -> (removed)
$ws.Cells.Item(113, 8).Value = 1407.1428  # H113: 1368.1818 -> 1407.1428
$ws.Cells.Item(113, 10).Value = 3000  # J113: 1640 -> 3000
$ws.Cells.Item(113, 12).Value = 3000  # L113: 1640 -> 3000
$ws.Cells.Item(113, 14).Value = -7340  # N113: -5980 -> -7340
$ws.Cells.Item(138, 8).Value = 46156.668  # H138: 49184.5 -> 46156.668
$ws.Cells.Item(138, 10).Value = 46156.668  # J138: 49184.5 -> 46156.668
$ws.Cells.Item(138, 12).Value = 46156.668  # L138: 49184.5 -> 46156.668
$ws.Cells.Item(138, 14).Value = -56436.668  # N138: -59464.5 -> -56436.668
$ws.Cells.Item(141, 8).Value = 53406  # H141: 57510 -> 53406
$ws.Cells.Item(141, 10).Value = 53406  # J141: 57510 -> 53406
$ws.Cells.Item(141, 12).Value = 53406  # L141: 57510 -> 53406
$ws.Cells.Item(141, 14).Value = -63766  # N141: -67870 -> -63766

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(22, 8).Value = 3262.5  # H22: 1475 -> 3262.5
$ws.Cells.Item(22, 9).Value = 350  # I22: 1633.3334 -> 350
$ws.Cells.Item(22, 10).Value = 12000  # J22: 1000 -> 12000
$ws.Cells.Item(22, 11).Value = 350  # K22: 1633.3334 -> 350
$ws.Cells.Item(22, 12).Value = 12000  # L22: 1000 -> 12000
$ws.Cells.Item(22, 13).Value = -55  # M22: -1338.3334 -> -55
$ws.Cells.Item(22, 14).Value = -12590  # N22: -1590 -> -12590
$ws.Cells.Item(27, 8).Value = 3262.5  # H27: 1475 -> 3262.5
$ws.Cells.Item(27, 9).Value = 350  # I27: 1633.3334 -> 350
$ws.Cells.Item(27, 10).Value = 12000  # J27: 1000 -> 12000
$ws.Cells.Item(27, 11).Value = 350  # K27: 1633.3334 -> 350
$ws.Cells.Item(27, 12).Value = 12000  # L27: 1000 -> 12000
$ws.Cells.Item(27, 13).Value = -243  # M27: -1526.3334 -> -243
$ws.Cells.Item(27, 14).Value = -12214  # N27: -1214 -> -12214
$ws.Cells.Item(61, 8).Value = 8554.467000000001  # H61: 11447.682 -> 8554.467000000001
$ws.Cells.Item(61, 9).Value = 11842.45  # I61: 13904.353 -> 11842.45
$ws.Cells.Item(61, 10).Value = 1978.5  # J61: 3095 -> 1978.5
$ws.Cells.Item(61, 11).Value = 11842.45  # K61: 13904.353 -> 11842.45
$ws.Cells.Item(61, 12).Value = 1978.5  # L61: 3095 -> 1978.5
$ws.Cells.Item(61, 13).Value = -11640.45  # M61: -13702.353 -> -11640.45
$ws.Cells.Item(61, 14).Value = -2382.5  # N61: -3499 -> -2382.5
$ws.Cells.Item(113, 8).Value = 8554.467000000001  # H113: 11447.682 -> 8554.467000000001
$ws.Cells.Item(113, 9).Value = 11842.45  # I113: 13904.353 -> 11842.45
$ws.Cells.Item(113, 10).Value = 1978.5  # J113: 3095 -> 1978.5
$ws.Cells.Item(113, 11).Value = 11842.45  # K113: 13904.353 -> 11842.45
$ws.Cells.Item(113, 12).Value = 1978.5  # L113: 3095 -> 1978.5
$ws.Cells.Item(113, 13).Value = -9672.450000000001  # M113: -11734.353 -> -9672.450000000001
$ws.Cells.Item(113, 14).Value = -6318.5  # N113: -7435 -> -6318.5
$ws.Cells.Item(132, 8).Value = 4571.163  # H132: 5014.5815 -> 4571.163
$ws.Cells.Item(132, 9).Value = 4679.95  # I132: 5128.5557 -> 4679.95
$ws.Cells.Item(132, 10).Value = 4087.6667  # J132: 4428.4287 -> 4087.6667
$ws.Cells.Item(132, 11).Value = 14039.85  # K132: 15385.6671 -> 14039.85
$ws.Cells.Item(132, 12).Value = 12263.0001  # L132: 13285.2861 -> 12263.0001
$ws.Cells.Item(132, 13).Value = -11509.85  # M132: -12855.6671 -> -11509.85
$ws.Cells.Item(132, 14).Value = -17323.0001  # N132: -18345.2861 -> -17323.0001
$ws.Cells.Item(139, 8).Value = 58000  # H139: 0 -> 58000
$ws.Cells.Item(139, 10).Value = 58000  # J139: 0 -> 58000
$ws.Cells.Item(139, 12).Value = 58000  # L139: 0 -> 58000
$ws.Cells.Item(139, 14).Value = -68280  # N139: None -> -68280

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(46, 8).Value = 50059.57  # H46: 49488.145 -> 50059.57
$ws.Cells.Item(46, 10).Value = 50059.57  # J46: 49488.145 -> 50059.57
$ws.Cells.Item(46, 12).Value = 50059.57  # L46: 49488.145 -> 50059.57
$ws.Cells.Item(46, 14).Value = -50521.57  # N46: -49950.145 -> -50521.57
$ws.Cells.Item(113, 8).Value = 836.4  # H113: 713.3077 -> 836.4
$ws.Cells.Item(113, 9).Value = 768  # I113: 1051 -> 768
$ws.Cells.Item(113, 10).Value = 865.7143  # J113: 651.9091 -> 865.7143
$ws.Cells.Item(113, 11).Value = 2304  # K113: 3153 -> 2304
$ws.Cells.Item(113, 12).Value = 2597.1429  # L113: 1955.7273 -> 2597.1429
$ws.Cells.Item(113, 13).Value = -134  # M113: -983 -> -134
$ws.Cells.Item(113, 14).Value = -6937.1429  # N113: -6295.7273 -> -6937.1429
$ws.Cells.Item(122, 8).Value = 35716130  # H122: 25001730 -> 35716130
$ws.Cells.Item(122, 9).Value = 125001000  # I122: 41667948 -> 125001000
$ws.Cells.Item(122, 10).Value = 2182  # J122: 2402.5 -> 2182
$ws.Cells.Item(122, 11).Value = 375003000  # K122: 125003844 -> 375003000
$ws.Cells.Item(122, 12).Value = 6546  # L122: 7207.5 -> 6546
$ws.Cells.Item(122, 13).Value = -375000550  # M122: -125001394 -> -375000550
$ws.Cells.Item(122, 14).Value = -11446  # N122: -12107.5 -> -11446
$ws.Cells.Item(134, 8).Value = 50059.57  # H134: 49488.145 -> 50059.57
$ws.Cells.Item(134, 10).Value = 50059.57  # J134: 49488.145 -> 50059.57
$ws.Cells.Item(134, 12).Value = 150178.71  # L134: 148464.435 -> 150178.71
$ws.Cells.Item(134, 14).Value = -155248.71  # N134: -153534.435 -> -155248.71
$ws.Cells.Item(137, 8).Value = 73045.836  # H137: 73989 -> 73045.836
$ws.Cells.Item(137, 10).Value = 73045.836  # J137: 73989 -> 73045.836
$ws.Cells.Item(137, 12).Value = 73045.836  # L137: 73989 -> 73045.836
$ws.Cells.Item(137, 14).Value = -83245.836  # N137: -84189 -> -83245.836
$ws.Cells.Item(139, 8).Value = 63367.5  # H139: 67473.75 -> 63367.5
$ws.Cells.Item(139, 10).Value = 63367.5  # J139: 67473.75 -> 63367.5
$ws.Cells.Item(139, 12).Value = 63367.5  # L139: 67473.75 -> 63367.5
$ws.Cells.Item(139, 14).Value = -73647.5  # N139: -77753.75 -> -73647.5
$ws.Cells.Item(141, 8).Value = 79000  # H141: 72470.55499999999 -> 79000
$ws.Cells.Item(141, 10).Value = 79000  # J141: 72470.55499999999 -> 79000
$ws.Cells.Item(141, 12).Value = 79000  # L141: 72470.55499999999 -> 79000
$ws.Cells.Item(141, 14).Value = -89360  # N141: -82830.55499999999 -> -89360
